# Update "想去人数" (want-to-go count) figures in column F across the
# refreshed data sheets (展览 / 演出 / 全部类型), matching the refreshed
# scrape snapshot generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 16519
$ws1.Range("F6").Value  = 15640
$ws1.Range("F9").Value  = 26
$ws1.Range("F10").Value = 109
$ws1.Range("F14").Value = 83
$ws1.Range("F15").Value = 1169
$ws1.Range("F19").Value = 563
$ws1.Range("F20").Value = 46
$ws1.Range("F21").Value = 48
$ws1.Range("F24").Value = 4
$ws1.Range("F26").Value = 281
$ws1.Range("F27").Value = 386
$ws1.Range("F29").Value = 123
$ws1.Range("F30").Value = 5840
$ws1.Range("F31").Value = 5271

# --- 演出 (Performance) sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 85

# --- 全部类型 (All types, merged listing) sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 16519
$ws4.Range("F6").Value  = 15640
$ws4.Range("F9").Value  = 26
$ws4.Range("F10").Value = 109
$ws4.Range("F14").Value = 83
$ws4.Range("F15").Value = 1169
$ws4.Range("F17").Value = 29
$ws4.Range("F19").Value = 563
$ws4.Range("F20").Value = 46
$ws4.Range("F21").Value = 48
$ws4.Range("F22").Value = 85
$ws4.Range("F26").Value = 4
$ws4.Range("F28").Value = 281
$ws4.Range("F29").Value = 386
$ws4.Range("F31").Value = 123
$ws4.Range("F32").Value = 5840
$ws4.Range("F34").Value = 5271
